$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N2").Value = 0.09250225208559804
$ws.Range("N3").Value = 0.001
$ws.Range("N4").Value = 0.00100000000000001
$ws.Range("N5").Value = 0.001000000000000019
$ws.Range("N7").Value = 0.15
$ws.Range("N9").Value = 0.001
$ws.Range("N10").Value = 0.15
$ws.Range("N11").Value = 0.04855061296321992
$ws.Range("N12").Value = 0.15
$ws.Range("N13").Value = 0.07286695165057186
$ws.Range("N14").Value = 0.1173398567334871
$ws.Range("N16").Value = 0.001000000000000002
$ws.Range("N19").Value = 0.1318109578632803
$ws.Range("N20").Value = 0.007632901877148632
$ws.Range("N21").Value = 0.06929646682669441
$ws.Range("B22").Value = 0.02147520643930769
$ws.Range("C22").Value = 0.04935615001843645
$ws.Range("D22").Value = 0.009755284474596169
$ws.Range("E22").Value = 0.03755228859330272
$ws.Range("F22").Value = 0.03426440281842462
$ws.Range("G22").Value = 0.02132798309349331
$ws.Range("H22").Value = 0.01300930426063261
$ws.Range("I22").Value = -0.01444327698244592
$ws.Range("J22").Value = 0.00256903922549441
$ws.Range("K22").Value = -0.01404324652842785
$ws.Range("L22").Value = 0.028263795820348
$ws.Range("M22").Value = 0.03185293027162392
$ws.Range("B23").Value = 1.021707458257471
$ws.Range("C23").Value = 1.050594453351678
$ws.Range("D23").Value = 1.009803022368154
$ws.Range("E23").Value = 1.038266285146013
$ws.Range("F23").Value = 1.034858189980533
$ws.Range("G23").Value = 1.021557050139088
$ws.Range("H23").Value = 1.013094293409309
$ws.Range("I23").Value = 0.9856605267860938
$ws.Range("J23").Value = 1.002572342034509
$ws.Range("K23").Value = 0.9860548998894985
$ws.Range("L23").Value = 1.028667006689906
$ws.Range("M23").Value = 1.032365664401943
$ws.Range("N23").Value = 1.127199742769974
